# This script rewrites the LR-pair data rows (rows 2-16) of Sheet1 to reflect
# the updated TPM-based NATMI computation: a new "Inflammatory-Mac" (D column)
# target-cluster block of rows is inserted between the existing "FAPs" and
# "Resolving-Mac" blocks for each sending cluster, and all numeric columns are
# refreshed with the new TPM-derived values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @{ Row=2; A="ECs"; B="Il18"; C="Il18r1"; D="FAPs"; E=1; F=0.3333333333333333; G=0.1611213333333333; H=0.483364; I=0.001886845364621106; J=0.001927375876930126; K=2; L=0.6666666666666666; M=0.3360656666666667; N=1.008197; O=0.4847958677304446; P=0.4847958677304446; Q=0.05414734830088889; R=0.487326134708; S=0.0009147348358146559; T=0.0009343838606990669 },
  @{ Row=3; A="ECs"; B="Il18"; C="Il18r1"; D="Inflammatory-Mac"; E=1; F=0.3333333333333333; G=0.1611213333333333; H=0.483364; I=0.001886845364621106; J=0.001927375876930126; K=2; L=0.6666666666666666; M=0.204299; N=0.612897; O=0.2947141609669403; P=0.2947141609669403; Q=0.03291692727866667; R=0.296252345508; S=0.0005560800485086697; T=0.0005680249644373828 },
  @{ Row=4; A="ECs"; B="Il18"; C="Il18r1"; D="Resolving-Mac"; E=1; F=0.3333333333333333; G=0.1611213333333333; H=0.483364; I=0.001886845364621106; J=0.001927375876930126; K=1; L=0.3333333333333333; M=0.152846; N=0.458538; O=0.2204899713026151; P=0.2204899713026151; Q=0.02462675131466667; R=0.221640761832; S=0.0004160304802977799; T=0.000424967051793676 },
  @{ Row=5; A="FAPs"; B="Il18"; C="Il18r1"; D="FAPs"; E=3; F=1; G=6.635478333333334; H=19.906435; I=0.07770616886214393; J=0.07937534159490063; K=2; L=0.6666666666666666; M=0.3360656666666667; N=1.008197; O=0.4847958677304446; P=0.4847958677304446; Q=2.229956449743889; R=20.069608047695; S=0.03767162956153152; T=0.0384808376049003 },
  @{ Row=6; A="FAPs"; B="Il18"; C="Il18r1"; D="Inflammatory-Mac"; E=3; F=1; G=6.635478333333334; H=19.906435; I=0.07770616886214393; J=0.07937534159490063; K=2; L=0.6666666666666666; M=0.204299; N=0.612897; O=0.2947141609669403; P=0.2947141609669403; Q=1.355621588021667; R=12.200594292195; S=0.02290110835816213; T=0.02339303719960542 },
  @{ Row=7; A="FAPs"; B="Il18"; C="Il18r1"; D="Resolving-Mac"; E=3; F=1; G=6.635478333333334; H=19.906435; I=0.07770616886214393; J=0.07937534159490063; K=1; L=0.3333333333333333; M=0.152846; N=0.458538; O=0.2204899713026151; P=0.2204899713026151; Q=1.014206321336667; R=9.127856892030001; S=0.01713343094245028; T=0.01750146679039491 },
  @{ Row=8; A="Inflammatory-Mac"; B="Il18"; C="Il18r1"; D="FAPs"; E=3; F=1; G=40.78183766666667; H=122.345513; I=0.4775843134495767; J=0.4878430963142499; K=2; L=0.6666666666666666; M=0.3360656666666667; N=1.008197; O=0.4847958677304446; P=0.4847958677304446; Q=13.70537546334011; R=123.348379170061; S=0.2315309016532362; T=0.2365043171939736 },
  @{ Row=9; A="Inflammatory-Mac"; B="Il18"; C="Il18r1"; D="Inflammatory-Mac"; E=3; F=1; G=40.78183766666667; H=122.345513; I=0.4775843134495767; J=0.4878430963142499; K=2; L=0.6666666666666666; M=0.204299; N=0.612897; O=0.2947141609669403; P=0.2947141609669403; Q=8.331688653462335; R=74.985197881161; S=0.1407508602292642; T=0.1437742688137684 },
  @{ Row=10; A="Inflammatory-Mac"; B="Il18"; C="Il18r1"; D="Resolving-Mac"; E=3; F=1; G=40.78183766666667; H=122.345513; I=0.4775843134495767; J=0.4878430963142499; K=1; L=0.3333333333333333; M=0.152846; N=0.458538; O=0.2204899713026151; P=0.2204899713026151; Q=6.233340759999334; R=56.100066839994; S=0.1053025515670763; T=0.1075645103065078 },
  @{ Row=11; A="MuSCs"; B="Il18"; C="Il18r1"; D="FAPs"; E=2; F=1; G=5.387083000000001; H=10.774166; I=0.06308657194606442; J=0.04296113827765565; K=2; L=0.6666666666666666; M=0.3360656666666667; N=1.008197; O=0.4847958677304446; P=0.4847958677304446; Q=1.810413639783667; R=10.862481838702; S=0.03058410938873142; T=0.02082738231000369 },
  @{ Row=12; A="MuSCs"; B="Il18"; C="Il18r1"; D="Inflammatory-Mac"; E=2; F=1; G=5.387083000000001; H=10.774166; I=0.06308657194606442; J=0.04296113827765565; K=2; L=0.6666666666666666; M=0.204299; N=0.612897; O=0.2947141609669403; P=0.2947141609669403; Q=1.100575669817; R=6.603454018902001; S=0.01859250611936489; T=0.01266125582168399 },
  @{ Row=13; A="MuSCs"; B="Il18"; C="Il18r1"; D="Resolving-Mac"; E=2; F=1; G=5.387083000000001; H=10.774166; I=0.06308657194606442; J=0.04296113827765565; K=1; L=0.3333333333333333; M=0.152846; N=0.458538; O=0.2204899713026151; P=0.2204899713026151; Q=0.8233940882180001; R=4.940364529308001; S=0.01390995643796811; T=0.009472500145967972 },
  @{ Row=14; A="Resolving-Mac"; B="Il18"; C="Il18r1"; D="FAPs"; E=3; F=1; G=32.42639166666667; H=97.279175; I=0.379736100377594; J=0.3878930479362637; K=2; L=0.6666666666666666; M=0.3360656666666667; N=1.008197; O=0.4847958677304446; P=0.4847958677304446; Q=10.89739693305278; R=98.07657239747499; S=0.1840944922911309; T=0.1880489467608679 },
  @{ Row=15; A="Resolving-Mac"; B="Il18"; C="Il18r1"; D="Inflammatory-Mac"; E=3; F=1; G=32.42639166666667; H=97.279175; I=0.379736100377594; J=0.3878930479362637; K=2; L=0.6666666666666666; M=0.204299; N=0.612897; O=0.2947141609669403; P=0.2947141609669403; Q=6.624679391108334; R=59.622114519975; S=0.1119136062116404; T=0.1143175741674451 },
  @{ Row=16; A="Resolving-Mac"; B="Il18"; C="Il18r1"; D="Resolving-Mac"; E=3; F=1; G=32.42639166666667; H=97.279175; I=0.379736100377594; J=0.3878930479362637; K=1; L=0.3333333333333333; M=0.152846; N=0.458538; O=0.2204899713026151; P=0.2204899713026151; Q=4.956244260683333; R=44.60619834614999; S=0.08372800187482267; T=0.08552652700795066 }
)

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 7).Value = $item.G
    $ws.Cells.Item($r, 8).Value = $item.H
    $ws.Cells.Item($r, 9).Value = $item.I
    $ws.Cells.Item($r, 10).Value = $item.J
    $ws.Cells.Item($r, 11).Value = $item.K
    $ws.Cells.Item($r, 12).Value = $item.L
    $ws.Cells.Item($r, 13).Value = $item.M
    $ws.Cells.Item($r, 14).Value = $item.N
    $ws.Cells.Item($r, 15).Value = $item.O
    $ws.Cells.Item($r, 16).Value = $item.P
    $ws.Cells.Item($r, 17).Value = $item.Q
    $ws.Cells.Item($r, 18).Value = $item.R
    $ws.Cells.Item($r, 19).Value = $item.S
    $ws.Cells.Item($r, 20).Value = $item.T
}
